$wb = $excel.ActiveWorkbook

# Sheet1 row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 3406.6584
$ws.Range("I138").Value = 4040.2856
$ws.Range("J138").Value = 3276.2058
$ws.Range("K138").Value = 12120.8568
$ws.Range("L138").Value = 9828.617400000001
$ws.Range("M138").Value = -6980.856800000001
$ws.Range("N138").Value = -20108.6174

# Sheet2 row 14
$ws = $wb.Worksheets.Item(2)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# Sheet2 row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 18725.807
$ws.Range("I32").Value = 3211.723
$ws.Range("K32").Value = 3211.723
$ws.Range("M32").Value = -2924.723

# Sheet2 row 43
$ws = $wb.Worksheets.Item(2)
$ws.Range("H43").Value = 8197.4
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 9496.75
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 9496.75
$ws.Range("M43").Value = -2687
$ws.Range("N43").Value = -10122.75

# Sheet2 row 97
$ws = $wb.Worksheets.Item(2)
$ws.Range("H97").Value = 45397.695
$ws.Range("I97").Value = 67969.07
$ws.Range("J97").Value = 3076.375
$ws.Range("K97").Value = 67969.07
$ws.Range("L97").Value = 3076.375
$ws.Range("M97").Value = -67473.07
$ws.Range("N97").Value = -4068.375

# Sheet2 row 110
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 20880424
$ws.Range("I110").Value = 26374400
$ws.Range("J110").Value = 3309.8
$ws.Range("K110").Value = 26374400
$ws.Range("L110").Value = 3309.8
$ws.Range("M110").Value = -26372355
$ws.Range("N110").Value = -7399.8

# Sheet2 row 124
$ws = $wb.Worksheets.Item(2)
$ws.Range("H124").Value = 21952.666
$ws.Range("J124").Value = 21952.666
$ws.Range("L124").Value = 21952.666
$ws.Range("N124").Value = -31772.666

# Sheet3 row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 27814038
$ws.Range("I107").Value = 41718292
$ws.Range("J107").Value = 5524.75
$ws.Range("K107").Value = 41718292
$ws.Range("L107").Value = 5524.75
$ws.Range("M107").Value = -41716372
$ws.Range("N107").Value = -9364.75

# Sheet4 row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 1439.2538
$ws.Range("I31").Value = 985.7619
$ws.Range("J31").Value = 2201.12
$ws.Range("K31").Value = 985.7619
$ws.Range("L31").Value = 2201.12
$ws.Range("M31").Value = -690.7619
$ws.Range("N31").Value = -2791.12

# Sheet4 row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 1439.2538
$ws.Range("I34").Value = 985.7619
$ws.Range("J34").Value = 2201.12
$ws.Range("K34").Value = 985.7619
$ws.Range("L34").Value = 2201.12
$ws.Range("M34").Value = -783.7619
$ws.Range("N34").Value = -2605.12

# Sheet4 row 124
$ws = $wb.Worksheets.Item(4)
$ws.Range("H124").Value = 41986.668
$ws.Range("J124").Value = 41986.668
$ws.Range("L124").Value = 41986.668
$ws.Range("N124").Value = -46896.668

# Sheet4 row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 1351.4706
$ws.Range("I134").Value = 1008.13794
$ws.Range("K134").Value = 3024.41382
$ws.Range("M134").Value = -489.4138199999998

# Sheet5 row 64
$ws = $wb.Worksheets.Item(5)
$ws.Range("H64").Value = 2333.1667
$ws.Range("I64").Value = 999.5
$ws.Range("K64").Value = 2998.5
$ws.Range("M64").Value = -2728.5

# Sheet5 row 67
$ws = $wb.Worksheets.Item(5)
$ws.Range("H67").Value = 2333.1667
$ws.Range("I67").Value = 999.5
$ws.Range("K67").Value = 2998.5
$ws.Range("M67").Value = -2062.5

# Sheet5 row 68
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 18919.86
$ws.Range("J68").Value = 27187.309
$ws.Range("L68").Value = 81561.927
$ws.Range("N68").Value = -83183.927

# Sheet5 row 71
$ws = $wb.Worksheets.Item(5)
$ws.Range("H71").Value = 18919.86
$ws.Range("J71").Value = 27187.309
$ws.Range("L71").Value = 244685.781
$ws.Range("N71").Value = -252797.781

# Sheet5 row 81
$ws = $wb.Worksheets.Item(5)
$ws.Range("H81").Value = 18520174
$ws.Range("I81").Value = 938.1818
$ws.Range("J81").Value = 47621828
$ws.Range("K81").Value = 2814.5454
$ws.Range("L81").Value = 142865484
$ws.Range("M81").Value = -1691.5454
$ws.Range("N81").Value = -142867730

# Sheet5 row 84
$ws = $wb.Worksheets.Item(5)
$ws.Range("H84").Value = 18520174
$ws.Range("I84").Value = 938.1818
$ws.Range("J84").Value = 47621828
$ws.Range("K84").Value = 8443.636199999999
$ws.Range("L84").Value = 428596452
$ws.Range("M84").Value = -2827.636199999999
$ws.Range("N84").Value = -428607684

# Sheet5 row 86
$ws = $wb.Worksheets.Item(5)
$ws.Range("H86").Value = 1260.8334
$ws.Range("I86").Value = 1000.4
$ws.Range("J86").Value = 1446.8572
$ws.Range("K86").Value = 3001.2
$ws.Range("L86").Value = 4340.571599999999
$ws.Range("M86").Value = -1815.2
$ws.Range("N86").Value = -6712.571599999999

# Sheet5 row 87
$ws = $wb.Worksheets.Item(5)
$ws.Range("H87").Value = 13059.417
$ws.Range("I87").Value = 3382.6
$ws.Range("J87").Value = 19971.428
$ws.Range("K87").Value = 10147.8
$ws.Range("L87").Value = 59914.284
$ws.Range("M87").Value = -8899.8
$ws.Range("N87").Value = -62410.284

# Sheet5 row 89
$ws = $wb.Worksheets.Item(5)
$ws.Range("H89").Value = 1260.8334
$ws.Range("I89").Value = 1000.4
$ws.Range("J89").Value = 1446.8572
$ws.Range("K89").Value = 9003.6
$ws.Range("L89").Value = 13021.7148
$ws.Range("M89").Value = -3075.6
$ws.Range("N89").Value = -24877.7148

# Sheet5 row 90
$ws = $wb.Worksheets.Item(5)
$ws.Range("H90").Value = 13059.417
$ws.Range("I90").Value = 3382.6
$ws.Range("J90").Value = 19971.428
$ws.Range("K90").Value = 30443.4
$ws.Range("L90").Value = 179742.852
$ws.Range("M90").Value = -24203.4
$ws.Range("N90").Value = -192222.852

# Sheet5 row 126
$ws = $wb.Worksheets.Item(5)
$ws.Range("H126").Value = 2840
$ws.Range("I126").Value = 2965
$ws.Range("J126").Value = 2790
$ws.Range("K126").Value = 8895
$ws.Range("L126").Value = 8370
$ws.Range("M126").Value = -3955
$ws.Range("N126").Value = -18250

# Sheet5 row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 1258795.4
$ws.Range("J131").Value = 1258795.4
$ws.Range("L131").Value = 3776386.2
$ws.Range("N131").Value = -3786466.2

# Sheet5 row 137
$ws = $wb.Worksheets.Item(5)
$ws.Range("H137").Value = 19610820
$ws.Range("I137").Value = 2968
$ws.Range("K137").Value = 8904
$ws.Range("M137").Value = -3804

# Sheet6 row 70
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 104064.4
$ws.Range("I70").Value = 157791.84
$ws.Range("J70").Value = 4284.857
$ws.Range("K70").Value = 157791.84
$ws.Range("L70").Value = 4284.857
$ws.Range("M70").Value = -157521.84
$ws.Range("N70").Value = -4824.857

# Sheet6 row 73
$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 104064.4
$ws.Range("I73").Value = 157791.84
$ws.Range("J73").Value = 4284.857
$ws.Range("K73").Value = 157791.84
$ws.Range("L73").Value = 4284.857
$ws.Range("M73").Value = -156855.84
$ws.Range("N73").Value = -6156.857

# Sheet6 row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 5885650
$ws.Range("J126").Value = 11766799
$ws.Range("L126").Value = 35300397
$ws.Range("N126").Value = -35305337

# Sheet7 row 19
$ws = $wb.Worksheets.Item(7)
$ws.Range("H19").Value = 19333.334
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 26500
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 26500
$ws.Range("M19").Value = -4830
$ws.Range("N19").Value = -26840

# Sheet7 row 138
$ws = $wb.Worksheets.Item(7)
$ws.Range("H138").Value = 45130
$ws.Range("I138").Value = 10390
$ws.Range("J138").Value = 62500
$ws.Range("K138").Value = 10390
$ws.Range("L138").Value = 62500
$ws.Range("M138").Value = -5250
$ws.Range("N138").Value = -72780

# Sheet8 row 96
$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 76925040
$ws.Range("I96").Value = 142859200
$ws.Range("J96").Value = 1843.3334
$ws.Range("K96").Value = 142859200
$ws.Range("L96").Value = 1843.3334
$ws.Range("M96").Value = -142857827
$ws.Range("N96").Value = -4589.3334

# Sheet8 row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 1645.6
$ws.Range("I136").Value = 840
$ws.Range("K136").Value = 2520
$ws.Range("M136").Value = 30
